# Append new scraped listings to the "ランサーズ" sheet, shifting the
# previous rows 2-5 down to rows 6-9, and refresh the capture timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-02-16 07:00:38"

# New data for rows 2-9 (row1 header untouched).
# Columns: A=取得日時 B=タイトル C=カテゴリ D=価格 E=締切 F=URL G=優先度スコア H=スキル概要
$rows = @(
    @{ B = "製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)"; D = "300,000 円 ~ 500,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5460562"; G = 435; H = "🔥AI,Ai ◆ツール,開発" },
    @{ B = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"; D = "20,000 円 ~ 50,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5217096"; G = 243; H = "🔥API ◆ツール" },
    @{ B = "施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集"; D = "300,000 円 ~ 500,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5460563"; G = 220; H = "◆開発,システム開発 ◇管理" },
    @{ B = "【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,200円程度)"; D = "50,000 円 ~ 100,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5488168"; G = 213; H = "🔥API ◇管理" },
    @{ B = "【社内用】Temu APIを使って受注データ、在庫データ、注文ステータスを更新してほしい"; D = "50,000 円 ~ 100,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5492576"; G = 188; H = "🔥API" },
    @{ B = "【エンジニア募集】香水自販機制御システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5492441"; G = 125; H = "◆開発,システム開発" },
    @{ B = "地域情報サイト 店舗データ自動収集・一括管理システム構築"; D = "1,000,000 円 ~ 3,000,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5492383"; G = 85; H = "◇サイト" },
    @{ B = "【急募】ダウンロードスクリプト開発とBOXアップロード依頼"; D = "100,000 円 ~ 200,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5492631"; G = 68; H = "◆開発" }
)

# Remove existing hyperlinks first so they can be rebuilt cleanly for the
# whole F column (row deletion/reassignment would otherwise leave stale
# relationships behind).
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = "システム開発"
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = "期限情報なし"
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $data.F) | Out-Null
}

# Column width tweaks: B 51 -> 52, H 19 -> 16 characters.
$ws.Columns.Item(2).ColumnWidth = 51.17
$ws.Columns.Item(8).ColumnWidth = 15.17
